$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '59.624.44'
$ws.Range("E2").Value = '  -1.94%  '

# Row 3
$ws.Range("D3").Value = '2.302.29'
$ws.Range("E3").Value = '  -2.58%  '

# Row 4
$ws.Range("E4").Value = '  +0.35%  '

# Row 5
$cell = $ws.Range("D5")
$cell.NumberFormat = "@"
$cell.Value = '541.26'
$cell.Style = "Normal"
$ws.Range("E5").Value = '  -1.57%  '

# Row 6
$cell = $ws.Range("D6")
$cell.NumberFormat = "@"
$cell.Value = '128.17'
$cell.Style = "Normal"
$ws.Range("E6").Value = '  -4.54%  '

# Row 7
$ws.Range("E7").Value = '  +0.36%  '

# Row 8
$ws.Range("E8").Value = '  -3.80%  '

# Row 9
$ws.Range("D9").Value = '2.299.57'
$ws.Range("E9").Value = '  -2.68%  '

# Row 10
$ws.Range("E10").Value = '  -0.97%  '

# Row 11
$cell = $ws.Range("D11")
$cell.NumberFormat = "@"
$cell.Value = '5.53'
$cell.Style = "Normal"
$ws.Range("E11").Value = '  -0.26%  '

# Row 12
$ws.Range("E12").Value = '  -0.85%  '

# Row 13
$cell = $ws.Range("D13")
$cell.NumberFormat = "@"
$cell.Value = '0.331'
$cell.Style = "Normal"
$ws.Range("E13").Value = '  -2.30%  '

# Row 14
$ws.Range("D14").Value = '2.713.53'
$ws.Range("E14").Value = '  -1.93%  '

# Row 15
$ws.Range("B15").Value = 'WrappedBTC'
$ws.Range("C15").Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range("D15").Value = '59.589.86'
$ws.Range("E15").Value = '  -1.44%  '

# Row 16
$ws.Range("B16").Value = 'Avalanche'
$ws.Range("C16").Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$cell = $ws.Range("D16")
$cell.NumberFormat = "@"
$cell.Value = '23.04'
$cell.Style = "Normal"
$ws.Range("E16").Value = '  -5.37%  '

# Row 17
$ws.Range("E17").Value = '  -2.52%  '

# Row 18
$ws.Range("D18").Value = '2.300.31'
$ws.Range("E18").Value = '  -1.95%  '

# Row 19
$cell = $ws.Range("D19")
$cell.NumberFormat = "@"
$cell.Value = '10.39'
$cell.Style = "Normal"
$ws.Range("E19").Value = '  -3.65%  '

# Row 20
$ws.Range("E20").Value = '  -5.24%  '

# Row 21
$cell = $ws.Range("D21")
$cell.NumberFormat = "@"
$cell.Value = '310.02'
$cell.Style = "Normal"
$ws.Range("E21").Value = '  -2.79%  '

# Row 22
$cell = $ws.Range("D22")
$cell.NumberFormat = "@"
$cell.Value = '6.49'
$cell.Style = "Normal"
$ws.Range("E22").Value = '  -5.89%  '

# Row 23
$cell = $ws.Range("D23")
$cell.NumberFormat = "@"
$cell.Value = '1.00'
$cell.Style = "Normal"
$ws.Range("E23").Value = '  -0.66%  '

# Row 24
$cell = $ws.Range("D24")
$cell.NumberFormat = "@"
$cell.Value = '63.00'
$cell.Style = "Normal"
$ws.Range("E24").Value = '  -0.55%  '

# Row 25
$ws.Range("E25").Value = '  -3.73%  '

# Row 26
$cell = $ws.Range("D26")
$cell.NumberFormat = "@"
$cell.Value = '0.998'
$cell.Style = "Normal"
$ws.Range("E26").Value = '  +0.03%  '

# Row 27
$ws.Range("E27").Value = '  -5.39%  '

# Row 28
$ws.Range("E28").Value = '  -1.12%  '

# Row 29
$ws.Range("E29").Value = '  +2.64%  '

# Row 30
$cell = $ws.Range("D30")
$cell.NumberFormat = "@"
$cell.Value = '172.02'
$cell.Style = "Normal"
$ws.Range("E30").Value = '  +0.14%  '

# Row 31
$ws.Range("E31").Value = '  -3.00%  '

# Row 32
$ws.Range("D32").Value = '0.0₃0715'
$ws.Range("E32").Value = '  -5.68%  '

# Row 33
$cell = $ws.Range("D33")
$cell.NumberFormat = "@"
$cell.Value = '5.78'
$cell.Style = "Normal"
$ws.Range("E33").Value = '  -2.91%  '

# Row 34
$ws.Range("E34").Value = '  -3.31%  '

# Row 35
$ws.Range("E35").Value = '  +0.04%  '

# Row 36
$ws.Range("E36").Value = '  -7.82%  '

# Row 37
$cell = $ws.Range("D37")
$cell.NumberFormat = "@"
$cell.Value = '17.65'
$cell.Style = "Normal"
$ws.Range("E37").Value = '  -2.48%  '

# Row 38
$ws.Range("E38").Value = '  +0.08%  '

# Row 39
$cell = $ws.Range("D39")
$cell.NumberFormat = "@"
$cell.Value = '3.98'
$cell.Style = "Normal"
$ws.Range("E39").Value = '  -4.71%  '

# Row 40
$cell = $ws.Range("D40")
$cell.NumberFormat = "@"
$cell.Value = '311.55'
$cell.Style = "Normal"
$ws.Range("E40").Value = '  -3.61%  '

# Row 41
$cell = $ws.Range("D41")
$cell.NumberFormat = "@"
$cell.Value = '37.57'
$cell.Style = "Normal"
$ws.Range("E41").Value = '  -2.17%  '

# Row 42
$ws.Range("E42").Value = '  -5.23%  '

# Row 43
$ws.Range("E43").Value = '  -6.43%  '

# Row 44
$ws.Range("E44").Value = '  -2.66%  '

# Row 45
$cell = $ws.Range("D45")
$cell.NumberFormat = "@"
$cell.Value = '0.0937'
$cell.Style = "Normal"
$ws.Range("E45").Value = '  -2.46%  '

# Row 46
$ws.Range("E46").Value = '  -0.13%  '

# Row 47
$ws.Range("B47").Value = 'BabyDogeCoin'
$ws.Range("C47").Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range("D47").Value = '0.0₆0230'
$ws.Range("E47").Value = '  +30.16%  '

# Row 48
$ws.Range("B48").Value = 'Hedera'
$ws.Range("C48").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$cell = $ws.Range("D48")
$cell.NumberFormat = "@"
$cell.Value = '0.0488'
$cell.Style = "Normal"
$ws.Range("E48").Value = '  -2.94%  '

# Row 49
$ws.Range("B49").Value = 'InjectiveProtocol'
$ws.Range("C49").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$cell = $ws.Range("D49")
$cell.NumberFormat = "@"
$cell.Value = '18.44'
$cell.Style = "Normal"
$ws.Range("E49").Value = '  -3.59%  '

# Row 50
$cell = $ws.Range("D50")
$cell.NumberFormat = "@"
$cell.Value = '0.0211'
$cell.Style = "Normal"
$ws.Range("E50").Value = '  -1.48%  '

# Row 51
$ws.Range("E51").Value = '  -0.27%  '
